$wb = $excel.ActiveWorkbook

# --- "Edit Repayment Schedule" sheet: leave behind a new selection (C8) ---
# (this sheet currently holds the tab-selected / active-cell state; once we
# select another sheet below, Excel drops tabSelected from this one but keeps
# whatever cell was last selected on it)
$wsEdit = $wb.Worksheets.Item("Edit Repayment Schedule")
$null = $wsEdit.Range("C8").Select()

# --- "Repayment schedule" sheet: insert a new (blank) "Variable Instalments"
# column before the existing "Late" column (currently column N) ---
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$null = $wsRepay.Columns("N").Insert()

# match the width Excel gives a freshly inserted column here (stored width 11)
$wsRepay.Columns("N").ColumnWidth = 10.166666666666666

# make "Repayment schedule" the active sheet/tab and leave the selection on K16
$null = $wsRepay.Select()
$null = $wsRepay.Range("K16").Select()
